$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Aman nyaman tepat waktu dan menarik banyak aktifitas musik dsb"
$ws.Range("A3").Value = "Cepat, tepat waktu dan bersih keretanya"
$ws.Range("A4").Value = "pelayanannya ramah ramah, terbaik, semoga bisa lebih baik"

$ws.Range("A4").Select()
